$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1115939503977046
$ws.Range("C2").Value = 0.8631666532838131
$ws.Range("B3").Value = 0.1226217426068196
$ws.Range("C3").Value = 0.9714525084005331
